# Apply updated cryptocurrency price/volume figures (refreshed data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain-text values (e.g. "1.000",
# "26.850.37", "  -0.89%  ") that must stay text -- mark them as Text format
# before assigning so Excel does not reinterpret them as numbers and silently
# drop significant trailing zeros / punctuation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.864.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5078"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07177"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8917"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.67"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.883.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07511"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.224"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008496"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.910.50"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.015"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.118.48"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.378"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.779"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.27%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.39"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.689"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.730"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09128"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7482"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.981"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.90%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.229"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.527"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5591"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01991"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.630"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.567"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.31%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4765"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.560"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.97"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.92%  "
